$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.346.07"
$ws.Range("E2").Value = "  -0.19%  "

$ws.Range("D3").Value = "2.642.24"
$ws.Range("E3").Value = "  -1.33%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "585.81"
$ws.Range("E5").Value = "  -3.01%  "

$ws.Range("D6").Value = "174.36"
$ws.Range("E6").Value = "  -2.72%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("E8").Value = "  -1.37%  "

$ws.Range("D9").Value = "2.640.06"
$ws.Range("E9").Value = "  -1.36%  "

$ws.Range("D10").Value = "0.171"
$ws.Range("E10").Value = "  -1.12%  "

$ws.Range("D11").Value = "0.169"
$ws.Range("E11").Value = "  +0.97%  "

$ws.Range("D12").Value = "0.355"
$ws.Range("E12").Value = "  -0.60%  "

$ws.Range("D13").Value = "4.94"
$ws.Range("E13").Value = "  -1.89%  "

$ws.Range("D14").Value = "3.127.89"
$ws.Range("E14").Value = "  -1.12%  "

$ws.Range("D15").Value = "0.0000185"
$ws.Range("E15").Value = "  -1.20%  "

$ws.Range("D16").Value = "72.210.49"
$ws.Range("E16").Value = "  -0.05%  "

$ws.Range("D17").Value = "25.87"
$ws.Range("E17").Value = "  -2.99%  "

$ws.Range("D18").Value = "2.660.35"
$ws.Range("E18").Value = "  -0.43%  "

$ws.Range("D19").Value = "12.13"
$ws.Range("E19").Value = "  +1.08%  "

$ws.Range("D20").Value = "7.92"
$ws.Range("E20").Value = "  -1.27%  "

$ws.Range("D21").Value = "373.85"
$ws.Range("E21").Value = "  -1.52%  "

$ws.Range("D22").Value = "4.14"
$ws.Range("E22").Value = "  -1.70%  "

$ws.Range("D23").Value = "2.06"
$ws.Range("E23").Value = "  -0.49%  "

$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("D25").Value = "70.94"
$ws.Range("E25").Value = "  -2.34%  "

$ws.Range("D26").Value = "4.24"
$ws.Range("E26").Value = "  -3.27%  "

$ws.Range("D27").Value = "9.58"
$ws.Range("E27").Value = "  -3.87%  "

$ws.Range("D28").Value = "2.778.10"
$ws.Range("E28").Value = "  -1.03%  "

$ws.Range("E29").Value = "  +0.41%  "

$ws.Range("D30").Value = "0.0₃0953"
$ws.Range("E30").Value = "  +0.00%  "

$ws.Range("D31").Value = "8.00"
$ws.Range("E31").Value = "  -2.11%  "

$ws.Range("D32").Value = "496.14"
$ws.Range("E32").Value = "  -5.30%  "

$ws.Range("E33").Value = "  -2.83%  "

$ws.Range("E34").Value = "  -1.52%  "

$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.07%  "

$ws.Range("D36").Value = "161.56"
$ws.Range("E36").Value = "  -1.95%  "

$ws.Range("D37").Value = "19.25"
$ws.Range("E37").Value = "  -2.03%  "

$ws.Range("E38").Value = "  +2.34%  "

$ws.Range("E39").Value = "  -1.33%  "

$ws.Range("D40").Value = "1.35"
$ws.Range("E40").Value = "  -2.79%  "

$ws.Range("E41").Value = "  -0.03%  "

$ws.Range("D42").Value = "1.73"
$ws.Range("E42").Value = "  -6.63%  "

$ws.Range("D43").Value = "2.58"
$ws.Range("E43").Value = "  -1.36%  "

$ws.Range("D44").Value = "4.90"
$ws.Range("E44").Value = "  -3.40%  "

$ws.Range("D45").Value = "0.326"
$ws.Range("E45").Value = "  -2.80%  "

$ws.Range("E46").Value = "  -0.68%  "

$ws.Range("D47").Value = "152.43"
$ws.Range("E47").Value = "  -0.70%  "

$ws.Range("D48").Value = "3.65"
$ws.Range("E48").Value = "  -2.84%  "

$ws.Range("D49").Value = "0.546"
$ws.Range("E49").Value = "  -1.21%  "

$ws.Range("D50").Value = "1.67"
$ws.Range("E50").Value = "  -3.21%  "

$ws.Range("D51").Value = "0.601"
$ws.Range("E51").Value = "  -0.55%  "
